$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fgf16"
$ws.Range("C2").Value = "Fgfr3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.229066
$ws.Range("H2").Value = 0.687198
$ws.Range("I2").Value = 0.1398528162861439
$ws.Range("J2").Value = 0.1398528162861439
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.166450999999999
$ws.Range("N2").Value = 12.499353
$ws.Range("O2").Value = 0.7126954333415383
$ws.Range("P2").Value = 0.7126954333415383
$ws.Range("Q2").Value = 0.9543922647659998
$ws.Range("R2").Value = 8.589530382893999
$ws.Range("S2").Value = 0.09967246350708789
$ws.Range("T2").Value = 0.09967246350708789
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fgf16"
$ws.Range("C3").Value = "Fgfr3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.229066
$ws.Range("H3").Value = 0.687198
$ws.Range("I3").Value = 0.1398528162861439
$ws.Range("J3").Value = 0.1398528162861439
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.6655859999999999
$ws.Range("N3").Value = 1.996758
$ws.Range("O3").Value = 0.1138523176430159
$ws.Range("P3").Value = 0.1138523176430159
$ws.Range("Q3").Value = 0.152463122676
$ws.Range("R3").Value = 1.372168104084
$ws.Range("S3").Value = 0.0159225672630804
$ws.Range("T3").Value = 0.0159225672630804
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fgf16"
$ws.Range("C4").Value = "Fgfr3"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.229066
$ws.Range("H4").Value = 0.687198
$ws.Range("I4").Value = 0.1398528162861439
$ws.Range("J4").Value = 0.1398528162861439
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.01401
$ws.Range("N4").Value = 3.04203
$ws.Range("O4").Value = 0.1734522490154458
$ws.Range("P4").Value = 0.1734522490154458
$ws.Range("Q4").Value = 0.23227521466
$ws.Range("R4").Value = 2.09047693194
$ws.Range("S4").Value = 0.02425778551597564
$ws.Range("T4").Value = 0.02425778551597563
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Fgf16"
$ws.Range("C5").Value = "Fgfr3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.408841666666667
$ws.Range("H5").Value = 4.226525
$ws.Range("I5").Value = 0.860147183713856
$ws.Range("J5").Value = 0.8601471837138561
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.166450999999999
$ws.Range("N5").Value = 12.499353
$ws.Range("O5").Value = 0.7126954333415383
$ws.Range("P5").Value = 0.7126954333415383
$ws.Range("Q5").Value = 5.869869770924999
$ws.Range("R5").Value = 52.82882793832499
$ws.Range("S5").Value = 0.6130229698344504
$ws.Range("T5").Value = 0.6130229698344505
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Fgf16"
$ws.Range("C6").Value = "Fgfr3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.408841666666667
$ws.Range("H6").Value = 4.226525
$ws.Range("I6").Value = 0.860147183713856
$ws.Range("J6").Value = 0.8601471837138561
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.6655859999999999
$ws.Range("N6").Value = 1.996758
$ws.Range("O6").Value = 0.1138523176430159
$ws.Range("P6").Value = 0.1138523176430159
$ws.Range("Q6").Value = 0.9377052895499998
$ws.Range("R6").Value = 8.43934760595
$ws.Range("S6").Value = 0.09792975037993545
$ws.Range("T6").Value = 0.09792975037993547
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Fgf16"
$ws.Range("C7").Value = "Fgfr3"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.408841666666667
$ws.Range("H7").Value = 4.226525
$ws.Range("I7").Value = 0.860147183713856
$ws.Range("J7").Value = 0.8601471837138561
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.01401
$ws.Range("N7").Value = 3.04203
$ws.Range("O7").Value = 0.1734522490154458
$ws.Range("P7").Value = 0.1734522490154458
$ws.Range("Q7").Value = 1.428579538416667
$ws.Range("R7").Value = 12.85721584575
$ws.Range("S7").Value = 0.1491944634994702
$ws.Range("T7").Value = 0.1491944634994702